$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 ---
$ws.Range("A22").Value = "dbx_v1_s50_r16_a90_b0_da0_de0_dr0"
$ws.Range("B22").Value = "NO"

# --- Row 23 ---
$ws.Range("A23").Value = "dbx_v1_s50_r16_am90_b0_da0_de0_dr0"
$ws.Range("B23").Value = "NO"

# Row heights / thick bottom border (mirrors existing rows' ht=15.75 thickBot=1 look)
$ws.Rows.Item(22).RowHeight = 15.75
$ws.Rows.Item(23).RowHeight = 15.75

# --- A22 / A23 formatting: same as A21 (wrapText, Arial 10, full medium-grey box border) ---
foreach ($addr in @("A22","A23")) {
  $c = $ws.Range($addr)
  $c.Font.Name = "Arial"
  $c.Font.Size = 10
  $c.Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleNone
  $c.WrapText = $true
  foreach ($edge in 7,8,9,10) {
    $b = $c.Borders.Item($edge)
    $b.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
    $b.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlMedium
    $b.Color = 13421772
  }
}

# --- B22: left + right medium-grey borders only, underlined Arial 10 ---
$b22 = $ws.Range("B22")
$b22.Font.Name = "Arial"
$b22.Font.Size = 10
$b22.Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleSingle
$b22.WrapText = $true
$b22.Borders.Item(7).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$b22.Borders.Item(7).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlMedium
$b22.Borders.Item(7).Color = 13421772
$b22.Borders.Item(10).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$b22.Borders.Item(10).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlMedium
$b22.Borders.Item(10).Color = 13421772
$b22.Borders.Item(8).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$b22.Borders.Item(9).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

# --- B23: left medium-grey border only, underlined Arial 10 ---
$b23 = $ws.Range("B23")
$b23.Font.Name = "Arial"
$b23.Font.Size = 10
$b23.Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleSingle
$b23.WrapText = $true
$b23.Borders.Item(7).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$b23.Borders.Item(7).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlMedium
$b23.Borders.Item(7).Color = 13421772
$b23.Borders.Item(8).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$b23.Borders.Item(9).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$b23.Borders.Item(10).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

# --- Selection: B23 ---
$ws.Range("B23").Select()
